# Scheduled runner update: refresh Universalis market-board price snapshots
# and recompute Leve profit columns (H:N) across the Kujata_Profits crafting sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3872.111
$ws.Range("I64").Value = 3848.1667
$ws.Range("K64").Value = 3848.1667
$ws.Range("M64").Value = -3600.1667
$ws.Range("H67").Value = 3872.111
$ws.Range("I67").Value = 3848.1667
$ws.Range("K67").Value = 3848.1667
$ws.Range("M67").Value = -2990.1667
$ws.Range("H69").Value = 3989.75
$ws.Range("J69").Value = 3989.75
$ws.Range("L69").Value = 11969.25
$ws.Range("N69").Value = -13717.25
$ws.Range("H70").Value = 1494.3334
$ws.Range("I70").Value = 1654
$ws.Range("J70").Value = 1334.6666
$ws.Range("K70").Value = 4962
$ws.Range("L70").Value = 4003.9998
$ws.Range("M70").Value = -4692
$ws.Range("N70").Value = -4543.9998
$ws.Range("H72").Value = 3989.75
$ws.Range("J72").Value = 3989.75
$ws.Range("L72").Value = 35907.75
$ws.Range("N72").Value = -44643.75
$ws.Range("H73").Value = 1494.3334
$ws.Range("I73").Value = 1654
$ws.Range("J73").Value = 1334.6666
$ws.Range("K73").Value = 4962
$ws.Range("L73").Value = 4003.9998
$ws.Range("M73").Value = -4026
$ws.Range("N73").Value = -5875.9998
$ws.Range("H80").Value = 615.08
$ws.Range("I80").Value = 545.2308
$ws.Range("J80").Value = 690.75
$ws.Range("K80").Value = 1635.6924
$ws.Range("L80").Value = 2072.25
$ws.Range("M80").Value = -637.6924000000001
$ws.Range("N80").Value = -4068.25
$ws.Range("H83").Value = 615.08
$ws.Range("I83").Value = 545.2308
$ws.Range("J83").Value = 690.75
$ws.Range("K83").Value = 4907.077200000001
$ws.Range("L83").Value = 6216.75
$ws.Range("M83").Value = 84.92279999999937
$ws.Range("N83").Value = -16200.75
$ws.Range("H100").Value = 1045
$ws.Range("I100").Value = 683.3333
$ws.Range("J100").Value = 2130
$ws.Range("K100").Value = 683.3333
$ws.Range("L100").Value = 2130
$ws.Range("M100").Value = -142.3333
$ws.Range("N100").Value = -3212
$ws.Range("H112").Value = 2478.6
$ws.Range("J112").Value = 2923.4375
$ws.Range("L112").Value = 8770.3125
$ws.Range("N112").Value = -10986.3125
$ws.Range("H137").Value = 1384.9219
$ws.Range("I137").Value = 1240.7
$ws.Range("K137").Value = 3722.1
$ws.Range("M137").Value = -1172.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H74").Value = 3296.6843
$ws.Range("I74").Value = 2818.5833
$ws.Range("K74").Value = 2818.5833
$ws.Range("M74").Value = -1944.5833
$ws.Range("H77").Value = 3296.6843
$ws.Range("I77").Value = 2818.5833
$ws.Range("K77").Value = 14092.9165
$ws.Range("M77").Value = -9724.916499999999
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("N96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("H132").Value = 3075.75
$ws.Range("I132").Value = 2360.1428
$ws.Range("K132").Value = 7080.428400000001
$ws.Range("M132").Value = -4550.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1369.7593
$ws.Range("I31").Value = 1359.3572
$ws.Range("J31").Value = 1406.1666
$ws.Range("K31").Value = 1359.3572
$ws.Range("L31").Value = 1406.1666
$ws.Range("M31").Value = -1064.3572
$ws.Range("N31").Value = -1996.1666
$ws.Range("H34").Value = 1369.7593
$ws.Range("I34").Value = 1359.3572
$ws.Range("J34").Value = 1406.1666
$ws.Range("K34").Value = 1359.3572
$ws.Range("L34").Value = 1406.1666
$ws.Range("M34").Value = -1157.3572
$ws.Range("N34").Value = -1810.1666
$ws.Range("H62").Value = 16672238
$ws.Range("I62").Value = 8235.714
$ws.Range("K62").Value = 8235.714
$ws.Range("M62").Value = -7611.714
$ws.Range("H65").Value = 16672238
$ws.Range("I65").Value = 8235.714
$ws.Range("K65").Value = 41178.57
$ws.Range("M65").Value = -38058.57
$ws.Range("H132").Value = 2647.5386
$ws.Range("I132").Value = 1940
$ws.Range("K132").Value = 5820
$ws.Range("M132").Value = -3290

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 3399.8333
$ws.Range("J81").Value = 3699.889
$ws.Range("L81").Value = 11099.667
$ws.Range("N81").Value = -13345.667
$ws.Range("H84").Value = 3399.8333
$ws.Range("J84").Value = 3699.889
$ws.Range("L84").Value = 33299.001
$ws.Range("N84").Value = -44531.001
$ws.Range("H87").Value = 1491.8572
$ws.Range("J87").Value = 1981
$ws.Range("L87").Value = 5943
$ws.Range("N87").Value = -8439
$ws.Range("H90").Value = 1491.8572
$ws.Range("J90").Value = 1981
$ws.Range("L90").Value = 17829
$ws.Range("N90").Value = -30309
$ws.Range("H98").Value = 1499
$ws.Range("J98").Value = 533.5
$ws.Range("L98").Value = 1600.5
$ws.Range("N98").Value = -4596.5
$ws.Range("H131").Value = 25003384
$ws.Range("I131").Value = 100000480
$ws.Range("J131").Value = 4352.967
$ws.Range("K131").Value = 300001440
$ws.Range("L131").Value = 13058.901
$ws.Range("M131").Value = -299996400
$ws.Range("N131").Value = -23138.901
$ws.Range("H138").Value = 2462.383
$ws.Range("I138").Value = 2517.0715
$ws.Range("J138").Value = 2439.182
$ws.Range("K138").Value = 7551.2145
$ws.Range("L138").Value = 7317.545999999999
$ws.Range("M138").Value = -2411.2145
$ws.Range("N138").Value = -17597.546
$ws.Range("H140").Value = 21682.32
$ws.Range("I140").Value = 42087.6
$ws.Range("K140").Value = 126262.8
$ws.Range("M140").Value = -121082.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 15000.125
$ws.Range("J63").Value = 15000.125
$ws.Range("L63").Value = 15000.125
$ws.Range("N63").Value = -16372.125
$ws.Range("H66").Value = 15000.125
$ws.Range("J66").Value = 15000.125
$ws.Range("L66").Value = 45000.375
$ws.Range("N66").Value = -51864.375
$ws.Range("H80").Value = 5722.222
$ws.Range("I80").Value = 6325
$ws.Range("J80").Value = 5240
$ws.Range("K80").Value = 6325
$ws.Range("L80").Value = 5240
$ws.Range("M80").Value = -5327
$ws.Range("N80").Value = -7236
$ws.Range("H83").Value = 5722.222
$ws.Range("I83").Value = 6325
$ws.Range("J83").Value = 5240
$ws.Range("K83").Value = 31625
$ws.Range("L83").Value = 26200
$ws.Range("M83").Value = -26633
$ws.Range("N83").Value = -36184

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 1760
$ws.Range("I4").Value = 800
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 800
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = -687
$ws.Range("N4").Value = -2226
$ws.Range("H28").Value = 1760
$ws.Range("I28").Value = 800
$ws.Range("J28").Value = 2000
$ws.Range("K28").Value = 800
$ws.Range("L28").Value = 2000
$ws.Range("M28").Value = -568
$ws.Range("N28").Value = -2464
$ws.Range("H37").Value = 1760
$ws.Range("I37").Value = 800
$ws.Range("J37").Value = 2000
$ws.Range("K37").Value = 800
$ws.Range("L37").Value = 2000
$ws.Range("M37").Value = -693
$ws.Range("N37").Value = -2214
$ws.Range("H55").Value = 241.75757
$ws.Range("I55").Value = 156.08696
$ws.Range("K55").Value = 156.08696
$ws.Range("M55").Value = 16.91304
$ws.Range("H93").Value = 879.05
$ws.Range("I93").Value = 841.1053000000001
$ws.Range("J93").Value = 1600
$ws.Range("K93").Value = 841.1053000000001
$ws.Range("L93").Value = 1600
$ws.Range("M93").Value = 406.8946999999999
$ws.Range("N93").Value = -4096

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 250001250
$ws.Range("J62").Value = 2500
$ws.Range("L62").Value = 2500
$ws.Range("N62").Value = -3748
$ws.Range("H65").Value = 250001250
$ws.Range("J65").Value = 2500
$ws.Range("L65").Value = 12500
$ws.Range("N65").Value = -18740
$ws.Range("H81").Value = 1309.5
$ws.Range("I81").Value = 1513.6666
$ws.Range("J81").Value = 1187
$ws.Range("K81").Value = 3027.3332
$ws.Range("L81").Value = 2374
$ws.Range("M81").Value = -1966.3332
$ws.Range("N81").Value = -4496
$ws.Range("H84").Value = 1309.5
$ws.Range("I84").Value = 1513.6666
$ws.Range("J84").Value = 1187
$ws.Range("K84").Value = 15136.666
$ws.Range("L84").Value = 11870
$ws.Range("M84").Value = -9832.666000000001
$ws.Range("N84").Value = -22478
$ws.Range("H122").Value = 11364927
$ws.Range("J122").Value = 1134.1666
$ws.Range("L122").Value = 3402.4998
$ws.Range("N122").Value = -8302.4998
$ws.Range("H135").Value = 98174.75
$ws.Range("J135").Value = 98174.75
$ws.Range("L135").Value = 98174.75
$ws.Range("N135").Value = -108314.75

